$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: add N4 value
$ws.Range("N4").Value = "Locations and Pallets"

# Row 5: N5 change text
$ws.Range("N5").Value = "Product Number, Pallets, From/To locations"

# Row 6: N6 change text
$ws.Range("N6").Value = "Product Number, Pallets"

# Update sheet view: scroll so column M is at the left edge, and select
# the full row 4 (activeCell A4, sqref A4:XFD4)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13 | Out-Null
$ws.Rows("4:4").Select() | Out-Null
